# edit.ps1 -- apply the "Add files via upload" revision to the ECG user
# manual.  The change:
#   1. The first bullet ("Open the file (ECGApp_final.m) in Matlab.")
#      becomes a new "Download" bullet listing the files to grab.
#   2. A new bullet with the original "Open the file (ECGApp_final.m) in
#      Matlab." text is inserted right after it.
#   3. A new bullet describing how to change the artifact width is
#      inserted right after the "Corrected (after removing pacing
#      artifacts) ..." bullet.
#
# We drive this with Range.InsertXML, which replaces the *content* of a
# Range with an OOXML fragment (wrapped in the usual single-part
# pkg:package/pkg:part/pkg:xmlData envelope Word uses for WordOpenXML),
# so every run / proofErr / bold marker ends up exactly as authored.

$d = $word.ActiveDocument

function New-WordPackageXml([string]$bodyXml) {
    return '<?xml version="1.0" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $bodyXml + '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
}

# Standard numbered-list paragraph properties used throughout this document.
$listPPr = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>'

# ---------------------------------------------------------------------
# 1) Paragraph 2 ("Open the file (ECGApp_final.m) in Matlab.") becomes
#    the new "Download ..." bullet.
# ---------------------------------------------------------------------
$downloadRuns =
    '<w:r><w:t xml:space="preserve">Download </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>ECGApp_final</w:t></w:r>' +
    '<w:r><w:t>.m</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve">, </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>kors.m</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>, Test_file1.mat, Test_file2.mat.</w:t></w:r>'

$p2 = $d.Paragraphs.Item(2)
$p2.Range.InsertXML((New-WordPackageXml ('<w:p>' + $listPPr + $downloadRuns + '</w:p>')))

# ---------------------------------------------------------------------
# 2) Insert a brand-new bullet right after it, holding the text that
#    used to live in paragraph 2 ("Open the file (ECGApp_final.m) in
#    Matlab.").
# ---------------------------------------------------------------------
$openFileRuns =
    '<w:r><w:t>Open the file (</w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>ECGApp_final</w:t></w:r>' +
    '<w:r><w:t>.m</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve">) in </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Matlab</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>.</w:t></w:r>'

$p2 = $d.Paragraphs.Item(2)
$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs.Item(3)
$p3.Range.InsertXML((New-WordPackageXml ('<w:p>' + $listPPr + $openFileRuns + '</w:p>')))

# ---------------------------------------------------------------------
# 3) Insert a new bullet about the artifact-width switch right after the
#    "Corrected (after removing pacing artifacts) ..." bullet.
# ---------------------------------------------------------------------
$bold = '<w:rPr><w:b/><w:bCs/></w:rPr>'
$widthRuns =
    '<w:r><w:t xml:space="preserve">To change the artifact width, enable the switch </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t xml:space="preserve">under </w:t></w:r>' +
    '<w:r><w:t>' + [char]0x2018 + '</w:t></w:r>' +
    '<w:r>' + $bold + '<w:t>%</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r>' + $bold + '<w:t>Envelope/Sample</w:t></w:r>' +
    '<w:r>' + $bold + '<w:t>' + [char]0x2019 + '</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">, to increase the width enter values below 0.3, to decrease the width enter values above 0.3 (default 0.3), alternatively, switch the method to </w:t></w:r>' +
    '<w:r>' + $bold + '<w:t>Sample</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> and enter the amount of data to remove, this value will be removed both before and after the artifact peak</w:t></w:r>' +
    '<w:r><w:t>.</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>'

$correctedPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text.StartsWith('Corrected (after removing pacing artifacts)')) {
        $correctedPara = $candidate
        break
    }
}

$correctedPara.Range.InsertParagraphAfter()
$newIndex = $correctedPara.Index + 1
$widthPara = $d.Paragraphs.Item($newIndex)
$widthPara.Range.InsertXML((New-WordPackageXml ('<w:p>' + $listPPr + $widthRuns + '</w:p>')))

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
